# Add a new "2022-Q4" quarter sheet (with its fund-holdings data) right after
# the "总计" summary sheet, and add a corresponding summary row to "总计".
#
# Net effect (matches the commit "feat: add 2022-Q4 data"):
#   - New sheet "2022-Q4" inserted at position 2 (all the other quarter sheets
#     shift one slot to the right; "2020-Q4" ends up last).
#   - "总计" gets a new row 2 for 2022-Q4 (持有数量(只)=5, 持有市值(亿元)=0.07),
#     with every following row's 日期 label shifted down one quarter and the
#     running index (column A) renumbered to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet by copying the existing "2022-Q3"
#    sheet (so header/styles/number formats all match its siblings) and
#    dropping the copy right before it -> it lands at position 2.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, [System.Reflection.Missing]::Value)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Extend the copied sheet (which only had 2 data rows) with 3 more rows so it
# can hold all 5 fund rows, matching the formatting of the existing data rows.
$q4.Range("A2:H2").Copy()
$q4.Range("A4:H4").PasteSpecial(-4122)
$q4.Range("A5:H5").PasteSpecial(-4122)
$q4.Range("A6:H6").PasteSpecial(-4122)

# The fund-code/name/size/position columns are stored as text even though
# several look numeric (e.g. "011685", "0.24") - force Text format before
# writing so they are not auto-coerced to numbers.
$q4.Range("B2:G6").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "011685"
$q4.Range("C2").Value = "创金合信先进装备股票A"
$q4.Range("D2").Value = "0.24"
$q4.Range("E2").Value = "80.29"
$q4.Range("F2").Value = "9.87"
$q4.Range("G2").Value = "0.0237"
$q4.Range("H2").Value = 1

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "519615"
$q4.Range("C3").Value = "银河君尚灵活配置混合I"
$q4.Range("D3").Value = "1.83"
$q4.Range("E3").Value = "38.98"
$q4.Range("F3").Value = "1.00"
$q4.Range("G3").Value = "0.0183"
$q4.Range("H3").Value = 3

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "011686"
$q4.Range("C4").Value = "创金合信先进装备股票C"
$q4.Range("D4").Value = "0.18"
$q4.Range("E4").Value = "80.29"
$q4.Range("F4").Value = "9.87"
$q4.Range("G4").Value = "0.0178"
$q4.Range("H4").Value = 1

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "519613"
$q4.Range("C5").Value = "银河君尚灵活配置混合A"
$q4.Range("D5").Value = "1.17"
$q4.Range("E5").Value = "38.98"
$q4.Range("F5").Value = "1.00"
$q4.Range("G5").Value = "0.0117"
$q4.Range("H5").Value = 3

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "519614"
$q4.Range("C6").Value = "银河君尚灵活配置混合C"
$q4.Range("D6").Value = "0.16"
$q4.Range("E6").Value = "38.98"
$q4.Range("F6").Value = "1.00"
$q4.Range("G6").Value = "0.0016"
$q4.Range("H6").Value = 3

# ---------------------------------------------------------------------------
# 2. "总计": insert a new row 2 for 2022-Q4, shifting the existing rows down
#    (they just move down one quarter each; the data values themselves don't
#    change, only their row position and the running index in column A).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Match the formatting of the row that got pushed down (same bold/bordered
# column-A style as all the other data rows).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.07000000000000001

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7
$total.Range("A10").Value = 8
